# Daily attendance processing - reverse the order of names/emails listed in
# the "Recorded By" column (column G) for every data row on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -ne $value -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ",\s*"

        $reversedParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $cell.Value2 = [string]::Join(", ", $reversedParts)
    }
}
